# Fixed a bug in Mask
# Re-applies the corrected row ordering/values for the reel-mask table on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{
    2  = @(902, 1, 0, 0, 0, 0)
    3  = @(601, 9, 60, 67, 60, 42)
    4  = @(201, 9, 30, 15, 45, 30)
    5  = @(1201, 2, 10, 10, 10, 10)
    6  = @(101, 9, 30, 15, 60, 15)
    7  = @(501, 9, 52, 30, 75, 45)
    9  = @(901, 16, 15, 45, 60, 60)
    10 = @(301, 6, 45, 30, 60, 45)
    11 = @(701, 3, 90, 45, 97, 15)
    12 = @(1203, 3, 15, 15, 15, 15)
    13 = @(1001, 18, 30, 75, 60, 72)
    14 = @(801, 3, 67, 65, 52, 45)
    16 = @(1, 0, 2, 2, 2, 2)
    17 = @(2, 0, 2, 2, 2, 2)
    18 = @(3, 0, 3, 3, 3, 3)
    19 = @(1101, 0, 15, 30, 30, 0)
    20 = @(502, 0, 4, 0, 0, 0)
    21 = @(802, 0, 4, 5, 4, 0)
    22 = @(402, 0, 0, 4, 0, 0)
    23 = @(602, 0, 0, 4, 0, 9)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
    }
}
